$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1238.2041
$ws.Range("I132").Value = 1323.8462
$ws.Range("J132").Value = 904.2
$ws.Range("K132").Value = 3971.5386
$ws.Range("L132").Value = 2712.6
$ws.Range("M132").Value = -1441.5386
$ws.Range("N132").Value = -7772.6
$ws.Range("H137").Value = 2753.2192
$ws.Range("I137").Value = 1998.362
$ws.Range("K137").Value = 5995.086
$ws.Range("M137").Value = -3445.086
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17119.98
$ws.Range("I32").Value = 2994.5974
$ws.Range("J32").Value = 64409.305
$ws.Range("K32").Value = 2994.5974
$ws.Range("L32").Value = 64409.305
$ws.Range("M32").Value = -2707.5974
$ws.Range("N32").Value = -64983.305
$ws.Range("H61").Value = 913.2712
$ws.Range("I61").Value = 792.3585
$ws.Range("J61").Value = 1981.3334
$ws.Range("K61").Value = 792.3585
$ws.Range("L61").Value = 1981.3334
$ws.Range("M61").Value = -580.3585
$ws.Range("N61").Value = -2405.3334
$ws.Range("H74").Value = 562.5769
$ws.Range("I74").Value = 345.3158
$ws.Range("J74").Value = 1152.2858
$ws.Range("K74").Value = 345.3158
$ws.Range("L74").Value = 1152.2858
$ws.Range("M74").Value = 528.6841999999999
$ws.Range("N74").Value = -2900.2858
$ws.Range("H77").Value = 562.5769
$ws.Range("I77").Value = 345.3158
$ws.Range("J77").Value = 1152.2858
$ws.Range("K77").Value = 1726.579
$ws.Range("L77").Value = 5761.429
$ws.Range("M77").Value = 2641.421
$ws.Range("N77").Value = -14497.429
$ws.Range("H132").Value = 1914.5671
$ws.Range("I132").Value = 1245.931
$ws.Range("J132").Value = 6223.5557
$ws.Range("K132").Value = 3737.793
$ws.Range("L132").Value = 18670.6671
$ws.Range("M132").Value = -1207.793
$ws.Range("N132").Value = -23730.6671
$ws.Range("H136").Value = 913.2712
$ws.Range("I136").Value = 792.3585
$ws.Range("J136").Value = 1981.3334
$ws.Range("K136").Value = 2377.0755
$ws.Range("L136").Value = 5944.0002
$ws.Range("M136").Value = 172.9245000000001
$ws.Range("N136").Value = -11044.0002
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 786.9231
$ws.Range("I134").Value = 711.64703
$ws.Range("J134").Value = 1298.8
$ws.Range("K134").Value = 2134.94109
$ws.Range("L134").Value = 3896.4
$ws.Range("M134").Value = 400.0589100000002
$ws.Range("N134").Value = -8966.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 763.0714
$ws.Range("I16").Value = 749.8889
$ws.Range("J16").Value = 786.8
$ws.Range("K16").Value = 749.8889
$ws.Range("L16").Value = 786.8
$ws.Range("M16").Value = -462.8889
$ws.Range("N16").Value = -1360.8
$ws.Range("H31").Value = 2797.5952
$ws.Range("I31").Value = 2578.258
$ws.Range("J31").Value = 3415.7273
$ws.Range("K31").Value = 2578.258
$ws.Range("L31").Value = 3415.7273
$ws.Range("M31").Value = -2283.258
$ws.Range("N31").Value = -4005.7273
$ws.Range("H34").Value = 2797.5952
$ws.Range("I34").Value = 2578.258
$ws.Range("J34").Value = 3415.7273
$ws.Range("K34").Value = 2578.258
$ws.Range("L34").Value = 3415.7273
$ws.Range("M34").Value = -2376.258
$ws.Range("N34").Value = -3819.7273
$ws.Range("H58").Value = 1443.45
$ws.Range("I58").Value = 1383.1072
$ws.Range("J58").Value = 1584.25
$ws.Range("K58").Value = 1383.1072
$ws.Range("L58").Value = 1584.25
$ws.Range("M58").Value = -1180.1072
$ws.Range("N58").Value = -1990.25
$ws.Range("H94").Value = 236448.77
$ws.Range("I94").Value = 667166.7
$ws.Range("J94").Value = 144152.08
$ws.Range("K94").Value = 667166.7
$ws.Range("L94").Value = 144152.08
$ws.Range("M94").Value = -666715.7
$ws.Range("N94").Value = -145054.08
$ws.Range("H99").Value = 1577.0769
$ws.Range("I99").Value = 1246.75
$ws.Range("J99").Value = 2105.6
$ws.Range("K99").Value = 1246.75
$ws.Range("L99").Value = 2105.6
$ws.Range("M99").Value = 251.25
$ws.Range("N99").Value = -5101.6
$ws.Range("H107").Value = 1416.0476
$ws.Range("I107").Value = 1573.2
$ws.Range("J107").Value = 1023.1667
$ws.Range("K107").Value = 1573.2
$ws.Range("L107").Value = 1023.1667
$ws.Range("M107").Value = 346.8
$ws.Range("N107").Value = -4863.1667
$ws.Range("H113").Value = 763.0714
$ws.Range("I113").Value = 749.8889
$ws.Range("J113").Value = 786.8
$ws.Range("K113").Value = 749.8889
$ws.Range("L113").Value = 786.8
$ws.Range("M113").Value = 1420.1111
$ws.Range("N113").Value = -5126.8
$ws.Range("H122").Value = 3344.8125
$ws.Range("I122").Value = 4065.182
$ws.Range("K122").Value = 12195.546
$ws.Range("M122").Value = -9745.545999999998
$ws.Range("H126").Value = 1577.0769
$ws.Range("I126").Value = 1246.75
$ws.Range("J126").Value = 2105.6
$ws.Range("K126").Value = 3740.25
$ws.Range("L126").Value = 6316.799999999999
$ws.Range("M126").Value = -1270.25
$ws.Range("N126").Value = -11256.8
$ws.Range("H132").Value = 1470.3793
$ws.Range("I132").Value = 1197.0571
$ws.Range("J132").Value = 1886.3043
$ws.Range("K132").Value = 3591.1713
$ws.Range("L132").Value = 5658.9129
$ws.Range("M132").Value = -1061.1713
$ws.Range("N132").Value = -10718.9129
$ws.Range("H134").Value = 1827.25
$ws.Range("I134").Value = 1843.96
$ws.Range("J134").Value = 1767.5714
$ws.Range("K134").Value = 5531.88
$ws.Range("L134").Value = 5302.7142
$ws.Range("M134").Value = -2996.88
$ws.Range("N134").Value = -10372.7142
$ws.Range("H136").Value = 1443.45
$ws.Range("I136").Value = 1383.1072
$ws.Range("J136").Value = 1584.25
$ws.Range("K136").Value = 4149.321599999999
$ws.Range("L136").Value = 4752.75
$ws.Range("M136").Value = -1599.321599999999
$ws.Range("N136").Value = -9852.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4417.3335
$ws.Range("I70").Value = 4270.8
$ws.Range("J70").Value = 5150
$ws.Range("K70").Value = 4270.8
$ws.Range("L70").Value = 5150
$ws.Range("M70").Value = -4000.8
$ws.Range("N70").Value = -5690
$ws.Range("H73").Value = 4417.3335
$ws.Range("I73").Value = 4270.8
$ws.Range("J73").Value = 5150
$ws.Range("K73").Value = 4270.8
$ws.Range("L73").Value = 5150
$ws.Range("M73").Value = -3334.8
$ws.Range("N73").Value = -7022
$ws.Range("H132").Value = 1760.2115
$ws.Range("I132").Value = 1493.3334
$ws.Range("J132").Value = 2881.1
$ws.Range("K132").Value = 4480.0002
$ws.Range("L132").Value = 8643.299999999999
$ws.Range("M132").Value = -1950.0002
$ws.Range("N132").Value = -13703.3
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1666.6842
$ws.Range("I7").Value = 1072.4445
$ws.Range("J7").Value = 2201.5
$ws.Range("K7").Value = 1072.4445
$ws.Range("L7").Value = 2201.5
$ws.Range("M7").Value = -960.4445000000001
$ws.Range("N7").Value = -2425.5
$ws.Range("H126").Value = 1666.6842
$ws.Range("I126").Value = 1072.4445
$ws.Range("J126").Value = 2201.5
$ws.Range("K126").Value = 3217.3335
$ws.Range("L126").Value = 6604.5
$ws.Range("M126").Value = -747.3335000000002
$ws.Range("N126").Value = -11544.5
$ws.Range("H132").Value = 2260.5715
$ws.Range("I132").Value = 2008.8727
$ws.Range("J132").Value = 3183.4666
$ws.Range("K132").Value = 6026.6181
$ws.Range("L132").Value = 9550.399800000001
$ws.Range("M132").Value = -3496.6181
$ws.Range("N132").Value = -14610.3998
